$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 2 (H) updates ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 550
$wsOff.Range("C2").Value = 383
$wsOff.Range("D2").Value = 130
$wsOff.Range("E2").Value = 58

# --- DEF sheet: row 2 (H) updates ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 523
$wsDef.Range("C2").Value = 362
$wsDef.Range("D2").Value = 115
$wsDef.Range("E2").Value = 42
$wsDef.Range("G2").Value = 9
